$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new person in row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "zxc"
$ws.Range("C7").Value = "f"
$ws.Range("D7").Value = 22

# Update the active selection to match the target state
$ws.Range("D7").Select()
